# Apply the content updates described by the diff to the "nakleyki_page_04" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (E2 / F2) - "Новый клиент - цена"
$ws.Range("E2").Value = "✅ Звісно! Друкуємо наліпки на якісних матеріалах.`n📄 Паперова самоклейка або 🌟 плівка (біла матова/прозора).`n💰 Прямокутні 90*50мм від 220 грн за 100 шт ⚡ 1-2 дні`n💰 Фігурні від 242 грн за 100 шт ⚡ 2-3 дні.`n❓ Чи маєте готовий для друку макет?"

$ws.Range("F2").Value = "✅ Конечно! Печатаем наклейки на качественных материалах.`n📄 Бумажная самоклейка или 🌟 пленка (белая матовая/прозрачная).`n💰 Прямоугольные 90х50мм от 220 грн за 100 шт ⚡ 1-2 дня`n💰 Фигурные от 242 грн за 100 шт ⚡ 2-3 дня`n❓ Есть ли у вас макет, готовый к печати?"

# Row 5 (E5) - "Другая типография" - drop the trailing tab character
$ws.Range("E5").Value = "Розуміємо, і це чудово! Ми дуже цінуємо роботу колег.`nПросто у нас такий підхід — ми відповідаємо за підсумкову якість наклейок на 100%, тому віддаємо перевагу векторним файлам.`nОсобливо важливо для фігурних наклейок — потрібна точність контурів для якісного вирізання.`nПереведемо у вектор від 💰 250 грн — і результат буде бездоганним."

# Row 6 (E6 / F6) - "Вопрос про материал"
$ws.Range("E6").Value = "📄 Друкуємо на двох типах матеріалів:`n🗞️ Самоклеючий папір — економний варіант, підходить для використання в приміщенні. Є варіант із посиленим клеєм, підходить для складних поверхонь і низьких температур.`n🌟 Плівка — біла матова або прозора, стійка до вологи та УФ-випромінювання, для зовнішнього використання.`nДля довготривалого застосування рекомендуємо плівку з ламінацією."

$ws.Range("F6").Value = "📄 Печатаем на двух типах материалов:`n🗞️ Бумажная самоклейка — экономичный вариант, подходит для внутреннего использования. Есть вариант с усиленным клеем, подходит для проблемных материалов и низких температурах.`n🌟 Пленка — белая матовая или прозрачная, стойкая к влаге и УФ, для наружного использования`nДля долговременного использования рекомендуем пленку с ламинацией."

# Row 8 (E8 / F8) - "Большие тиражи"
$ws.Range("E8").Value = "Чудово! Великі тиражі наклейок — наша сила.`nКвадратні 40х40мм: 5000 шт — від 2770 грн, 10000 шт — від 5185 грн.`nКруглі 40мм діаметром: 5000 шт — від 3330 грн, 10000 шт — від 6230 грн.`nРозкажіть точний тираж і форму — порахуємо найкращу ціну!        "

$ws.Range("F8").Value = "Отлично! Большие тиражи наклеек — наша сила.`nКвадратные 40х40: 5000 шт — от 2770 грн, 10000 шт — от 5185 грн.`nКруглые 40мм диаметр: 5000 шт — от 3330 грн, 10000 шт — от 6230 грн.`nРасскажите точный тираж и форму — посчитаем лучшую цену!"
